$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.504.43'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.912.91'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.31%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.54%  '

$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4790'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.50%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2845'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06702'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.68%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.44'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.03%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '102.87'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.97%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07797'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.30%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.922.64'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.95%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.205'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.83%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6710'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.29%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '275.39'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.48%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.558.33'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.98%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007480'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.65'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.43%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.391'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.38%  '

$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.309'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.41%  '

$ws.Range("E24").Value = '  -5.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '167.70'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.95%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.19'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.16%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.083'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.71%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.382'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.09971'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.44%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.576'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.45%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.516'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.63%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.256'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.78%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04736'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.21%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7267'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.112'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.41%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.720'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.57%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01910'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.90%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.620'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.322'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.46%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '73.86'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.26%  '

$ws.Range("E41").Value = '  -6.33%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '106.65'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8619'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.62%  '

$ws.Range("E44").Value = '  -4.00%  '

$ws.Range("E45").Value = '  +0.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.398'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.20%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '949.29'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.64%  '

$ws.Range("E48").Value = '  -3.62%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.66'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05802'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.56%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.701'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.51%  '

